$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# Append two new weather key/value rows (53 and 54), matching the
# existing key/value layout used throughout the sheet (column A = key,
# column B = localized value, with the same wrap-text style as column B
# elsewhere).
$ws.Range("A53").Value = "weather_haze"
$ws.Range("B53").Value = "Haze"
$ws.Range("B53").WrapText = $true

$ws.Range("A54").Value = "weather_dustStorm"
$ws.Range("B54").Value = "Dust Storm"
$ws.Range("B54").WrapText = $true

$ws.Range("A54").Select()
